$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 4
$ws.Range("F2").Value = 189
$ws.Range("H2").Value = "bedrooms"
$ws.Range("L2").Value = "stimuli/img_uxxo0.png"
$ws.Range("M2").Value = 71.74418604651163
$ws.Range("N2").Value = 48.44186046511628
$ws.Range("O2").Value = 60.09302325581395
$ws.Range("P2").Value = 43

# Row 3
$ws.Range("C3").Value = 4
$ws.Range("F3").Value = 190
$ws.Range("H3").Value = "bedrooms"
$ws.Range("L3").Value = "stimuli/img_2js6m.png"
$ws.Range("M3").Value = 40.02777777777778
$ws.Range("N3").Value = 20.88888888888889
$ws.Range("O3").Value = 30.45833333333334
$ws.Range("P3").Value = 36

# Row 4
$ws.Range("C4").Value = 4
$ws.Range("F4").Value = 191
$ws.Range("H4").Value = "bedrooms"
$ws.Range("L4").Value = "stimuli/img_zgg62.png"
$ws.Range("M4").Value = 82.18421052631579
$ws.Range("N4").Value = 63.52631578947368
$ws.Range("O4").Value = 72.85526315789474
$ws.Range("P4").Value = 38
$ws.Range("Q4").Value = 8
$ws.Range("R4").Value = 8
$ws.Range("S4").Value = 8

# Row 5
$ws.Range("C5").Value = 4
$ws.Range("F5").Value = 192
$ws.Range("H5").Value = "bedrooms"
$ws.Range("L5").Value = "stimuli/img_th7xh.png"
$ws.Range("M5").Value = 82.35897435897436
$ws.Range("N5").Value = 65.53846153846153
$ws.Range("O5").Value = 73.94871794871796
$ws.Range("P5").Value = 39

# Row 6
$ws.Range("C6").Value = 4
$ws.Range("F6").Value = 193
$ws.Range("H6").Value = "bedrooms"
$ws.Range("L6").Value = "stimuli/img_5yhyk.png"
$ws.Range("M6").Value = 46.375
$ws.Range("N6").Value = 31.325
$ws.Range("O6").Value = 38.85
$ws.Range("P6").Value = 40
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 2
$ws.Range("S6").Value = 2

# Row 7
$ws.Range("C7").Value = 4
$ws.Range("F7").Value = 194
$ws.Range("H7").Value = "bedrooms"
$ws.Range("L7").Value = "stimuli/img_i7vab.png"
$ws.Range("M7").Value = 86.40000000000001
$ws.Range("N7").Value = 67.8
$ws.Range("O7").Value = 77.09999999999999
$ws.Range("P7").Value = 35

# Row 8
$ws.Range("C8").Value = 4
$ws.Range("F8").Value = 195
$ws.Range("H8").Value = "bedrooms"
$ws.Range("L8").Value = "stimuli/img_h0hbk.png"
$ws.Range("M8").Value = 86.80952380952381
$ws.Range("N8").Value = 69.19047619047619
$ws.Range("O8").Value = 78
$ws.Range("P8").Value = 42

# Row 9
$ws.Range("C9").Value = 4
$ws.Range("F9").Value = 196
$ws.Range("H9").Value = "bedrooms"
$ws.Range("L9").Value = "stimuli/img_5m6x4.png"
$ws.Range("M9").Value = 80.23076923076923
$ws.Range("N9").Value = 58.41025641025641
$ws.Range("O9").Value = 69.32051282051282
$ws.Range("P9").Value = 39

# Row 10
$ws.Range("C10").Value = 4
$ws.Range("F10").Value = 197
$ws.Range("H10").Value = "bedrooms"
$ws.Range("L10").Value = "stimuli/img_le8uf.png"
$ws.Range("M10").Value = 12.88888888888889
$ws.Range("N10").Value = 9.222222222222221
$ws.Range("O10").Value = 11.05555555555556
$ws.Range("P10").Value = 36

# Row 11
$ws.Range("C11").Value = 4
$ws.Range("F11").Value = 198
$ws.Range("H11").Value = "bedrooms"
$ws.Range("L11").Value = "stimuli/img_0eflx.png"
$ws.Range("M11").Value = 76.05128205128206
$ws.Range("N11").Value = 53.53846153846154
$ws.Range("O11").Value = 64.7948717948718
$ws.Range("P11").Value = 39
$ws.Range("Q11").Value = 6
$ws.Range("R11").Value = 6
$ws.Range("S11").Value = 6

# Row 12
$ws.Range("C12").Value = 4
$ws.Range("F12").Value = 199
$ws.Range("H12").Value = "bedrooms"
$ws.Range("L12").Value = "stimuli/img_zv0dq.png"
$ws.Range("M12").Value = 76.86842105263158
$ws.Range("N12").Value = 52.71052631578947
$ws.Range("O12").Value = 64.78947368421052
$ws.Range("P12").Value = 38
$ws.Range("Q12").Value = 6
$ws.Range("R12").Value = 6
$ws.Range("S12").Value = 6

# Row 13
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 200
$ws.Range("H13").Value = "bedrooms"
$ws.Range("L13").Value = "stimuli/img_bklr1.png"
$ws.Range("M13").Value = 86.54761904761905
$ws.Range("N13").Value = 67.73809523809524
$ws.Range("O13").Value = 77.14285714285714
$ws.Range("P13").Value = 42
$ws.Range("Q13").Value = 9
$ws.Range("R13").Value = 9
$ws.Range("S13").Value = 9

# Row 14
$ws.Range("C14").Value = 4
$ws.Range("F14").Value = 201
$ws.Range("H14").Value = "bedrooms"
$ws.Range("L14").Value = "stimuli/img_qgbyn.png"
$ws.Range("M14").Value = 65.08108108108108
$ws.Range("N14").Value = 40.10810810810811
$ws.Range("O14").Value = 52.5945945945946
$ws.Range("P14").Value = 37

# Row 15
$ws.Range("C15").Value = 4
$ws.Range("F15").Value = 202
$ws.Range("H15").Value = "bedrooms"
$ws.Range("L15").Value = "stimuli/img_v8dra.png"
$ws.Range("M15").Value = 61.77272727272727
$ws.Range("N15").Value = 38.79545454545455
$ws.Range("O15").Value = 50.28409090909091
$ws.Range("P15").Value = 44

# Row 16
$ws.Range("C16").Value = 4
$ws.Range("F16").Value = 203

# Row 17
$ws.Range("C17").Value = 4
$ws.Range("F17").Value = 204
$ws.Range("H17").Value = "bedrooms"
$ws.Range("L17").Value = "stimuli/img_rvssl.png"
$ws.Range("M17").Value = 74.25
$ws.Range("N17").Value = 54.33333333333334
$ws.Range("O17").Value = 64.29166666666667
$ws.Range("P17").Value = 36

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("F18").Value = 205
$ws.Range("H18").Value = "bedrooms"
$ws.Range("L18").Value = "stimuli/img_x0u5z.png"
$ws.Range("M18").Value = 92
$ws.Range("N18").Value = 78.16216216216216
$ws.Range("O18").Value = 85.08108108108108
$ws.Range("P18").Value = 37

# Row 19
$ws.Range("C19").Value = 4
$ws.Range("F19").Value = 206
$ws.Range("H19").Value = "bedrooms"
$ws.Range("L19").Value = "stimuli/img_fqgem.png"
$ws.Range("M19").Value = 80.75
$ws.Range("N19").Value = 61.475
$ws.Range("O19").Value = 71.1125
$ws.Range("P19").Value = 40

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("F20").Value = 207
$ws.Range("H20").Value = "bedrooms"
$ws.Range("L20").Value = "stimuli/img_ybbmx.png"
$ws.Range("M20").Value = 55.24324324324324
$ws.Range("N20").Value = 36.75675675675676
$ws.Range("O20").Value = 46
$ws.Range("P20").Value = 37

# Row 21
$ws.Range("C21").Value = 4
$ws.Range("F21").Value = 208
$ws.Range("H21").Value = "bedrooms"
$ws.Range("L21").Value = "stimuli/img_oou46.png"
$ws.Range("M21").Value = 75.70270270270271
$ws.Range("N21").Value = 54.86486486486486
$ws.Range("O21").Value = 65.28378378378379
$ws.Range("P21").Value = 37

# Row 22
$ws.Range("C22").Value = 4
$ws.Range("F22").Value = 209
$ws.Range("H22").Value = "bedrooms"
$ws.Range("L22").Value = "stimuli/img_wyctg.png"
$ws.Range("M22").Value = 33.44736842105263
$ws.Range("N22").Value = 11.39473684210526
$ws.Range("O22").Value = 22.42105263157895
$ws.Range("P22").Value = 38

# Row 23
$ws.Range("C23").Value = 4
$ws.Range("F23").Value = 210
$ws.Range("H23").Value = "bedrooms"
$ws.Range("L23").Value = "stimuli/img_okvvw.png"
$ws.Range("M23").Value = 50.58333333333334
$ws.Range("N23").Value = 32.11111111111111
$ws.Range("O23").Value = 41.34722222222223
$ws.Range("P23").Value = 36

# Row 24
$ws.Range("C24").Value = 4
$ws.Range("F24").Value = 211
$ws.Range("H24").Value = "bedrooms"
$ws.Range("L24").Value = "stimuli/img_71mhq.png"
$ws.Range("M24").Value = 69.34210526315789
$ws.Range("N24").Value = 47.02631578947368
$ws.Range("O24").Value = 58.18421052631579
$ws.Range("P24").Value = 38

# Row 25
$ws.Range("C25").Value = 4
$ws.Range("F25").Value = 212
$ws.Range("H25").Value = "bedrooms"
$ws.Range("L25").Value = "stimuli/img_a9acb.png"
$ws.Range("M25").Value = 77.11428571428571
$ws.Range("N25").Value = 58.42857142857143
$ws.Range("O25").Value = 67.77142857142857
$ws.Range("P25").Value = 35

# Row 26
$ws.Range("C26").Value = 4
$ws.Range("F26").Value = 213
$ws.Range("H26").Value = "bedrooms"
$ws.Range("L26").Value = "stimuli/img_bj2gr.png"
$ws.Range("M26").Value = 65.25
$ws.Range("N26").Value = 44.8
$ws.Range("O26").Value = 55.025
$ws.Range("P26").Value = 40

# Row 27
$ws.Range("C27").Value = 4
$ws.Range("F27").Value = 214
$ws.Range("H27").Value = "bedrooms"
$ws.Range("L27").Value = "stimuli/img_2pk6v.png"
$ws.Range("M27").Value = 85.08108108108108
$ws.Range("N27").Value = 66.16216216216216
$ws.Range("O27").Value = 75.62162162162161
$ws.Range("P27").Value = 37

# Row 28
$ws.Range("C28").Value = 4
$ws.Range("F28").Value = 215
$ws.Range("H28").Value = "bedrooms"
$ws.Range("L28").Value = "stimuli/img_6ddrx.png"
$ws.Range("M28").Value = 82.2
$ws.Range("N28").Value = 63.68571428571428
$ws.Range("O28").Value = 72.94285714285715
$ws.Range("P28").Value = 35
$ws.Range("Q28").Value = 8
$ws.Range("R28").Value = 8
$ws.Range("S28").Value = 8

# Row 29
$ws.Range("C29").Value = 4
$ws.Range("F29").Value = 216
$ws.Range("H29").Value = "bedrooms"
$ws.Range("L29").Value = "stimuli/img_t2ioc.png"
$ws.Range("M29").Value = 88.18918918918919
$ws.Range("N29").Value = 74.05405405405405
$ws.Range("O29").Value = 81.12162162162161
$ws.Range("P29").Value = 37
